$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "xTestWordx"
$ws.Range("F2").Value = "Wassap"
$ws.Range("F3").Value = "Hahahaha"
$ws.Range("F4").Value = "Yowz"

$ws.Range("F5").Select()
